# projects_vineet.xlsx — "Add files via upload"
#
# The uploaded version of the sheet removes the project row for
# "Digitise Doctors' handwritten prescriptions using vision and deep
# learning techniques" / "BANGALORE- conditional" (previously row 8),
# which shifts every subsequent row up by one and drops the two shared
# strings that only that row used. It also adds a "T" marker in column D
# of the "Course5 intelligence customer segmentation" row (now row 5),
# and the active selection moves to B8 with the viewport scrolled so
# row 5 is at the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 8 ("Digitise Doctors' ..." / "BANGALORE- conditional").
# Excel shifts rows 9-11 up to become rows 8-10 automatically, and the two
# shared strings that were unique to that row are dropped from the workbook.
$ws.Rows.Item(8).Delete()

# The "Course5 intelligence customer segmentation" row (now row 5) gains a
# "T" value in column D, matching the other rows that already had it.
$ws.Range("D5").Value = "T"

# Update the view: scroll so row 5 is at the top and select B8.
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("B8").Select()
